$d = $word.ActiveDocument

# The last paragraph in the document currently contains only the
# "_GoBack" bookmark (no visible text). We need to:
#  1. Insert a brand new paragraph before it with the "3.try ..." text.
#  2. Insert a run with "  This solution only works for 10 and 100. "
#     text at the start of that (still) last paragraph, i.e. before the
#     bookmark.

$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range

# Step 1: insert a new paragraph before the bookmark paragraph containing
# the first sentence.
$lastRange.InsertBefore("3.try to predict the number she stops on by looking at even and odd numbers.`r")

# After the InsertBefore above, $d.Paragraphs.Last is still the bookmark
# paragraph (now preceded by the newly created one).
$bookmarkPara = $d.Paragraphs.Last
$bookmarkRange = $bookmarkPara.Range

# Step 2: insert text at the very beginning of the bookmark paragraph.
$insertPoint = $d.Range($bookmarkRange.Start, $bookmarkRange.Start)
$insertPoint.InsertBefore("  This solution only works for 10 and 100. ")
